$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 268; this shifts rows 268:282 down to 269:283
# and preserves all of their existing values/formatting.
$ws.Rows(268).Insert()

# Populate the newly inserted row 268 with a new weekly price record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the constant values used by the
# surrounding rows for this market/category/variety.
$ws.Range("A268").Value = 4
$ws.Range("B268").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C268").Value = "Los Lagos"
$ws.Range("D268").Value = 44746
$ws.Range("E268").Value = 10
$ws.Range("F268").Value = 100112003
$ws.Range("G268").Value = "Ajo"
$ws.Range("H268").Value = "Chino"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 80
$ws.Range("K268").Value = 22000
$ws.Range("L268").Value = 22000
$ws.Range("M268").Value = 22000
$ws.Range("N268").Value = "$/caja 10 kilos"
$ws.Range("O268").Value = "China"
$ws.Range("P268").Value = 2200
$ws.Range("Q268").Value = 10
$ws.Range("R268").Value = "Hortaliza"
